$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6728.9473
$ws.Range("I62").Value = 6132.615
$ws.Range("J62").Value = 8021
$ws.Range("K62").Value = 6132.615
$ws.Range("L62").Value = 8021
$ws.Range("M62").Value = -5508.615
$ws.Range("N62").Value = -9269
$ws.Range("H65").Value = 6728.9473
$ws.Range("I65").Value = 6132.615
$ws.Range("J65").Value = 8021
$ws.Range("K65").Value = 30663.075
$ws.Range("L65").Value = 40105
$ws.Range("M65").Value = -27543.075
$ws.Range("N65").Value = -46345
$ws.Range("H70").Value = 7603.5713
$ws.Range("I70").Value = 3023
$ws.Range("J70").Value = 10148.333
$ws.Range("K70").Value = 9069
$ws.Range("L70").Value = 30444.999
$ws.Range("M70").Value = -8799
$ws.Range("N70").Value = -30984.999
$ws.Range("H73").Value = 7603.5713
$ws.Range("I73").Value = 3023
$ws.Range("J73").Value = 10148.333
$ws.Range("K73").Value = 9069
$ws.Range("L73").Value = 30444.999
$ws.Range("M73").Value = -8133
$ws.Range("N73").Value = -32316.999
$ws.Range("H98").Value = 5410.2383
$ws.Range("I98").Value = 3070.5881
$ws.Range("K98").Value = 3070.5881
$ws.Range("M98").Value = -1572.5881
$ws.Range("H100").Value = 9499.75
$ws.Range("I100").Value = 7749.5
$ws.Range("K100").Value = 7749.5
$ws.Range("M100").Value = -7208.5
$ws.Range("H103").Value = 1575.381
$ws.Range("I103").Value = 1345.9
$ws.Range("K103").Value = 4037.7
$ws.Range("M103").Value = -3451.7
$ws.Range("H121").Value = 2228.3333
$ws.Range("J121").Value = 2228.3333
$ws.Range("L121").Value = 6684.999899999999
$ws.Range("N121").Value = -10178.9999
$ws.Range("H122").Value = 5410.2383
$ws.Range("I122").Value = 3070.5881
$ws.Range("K122").Value = 9211.764299999999
$ws.Range("M122").Value = -6761.764299999999
$ws.Range("H129").Value = 2193.6924
$ws.Range("J129").Value = 4991
$ws.Range("L129").Value = 14973
$ws.Range("N129").Value = -24973
$ws.Range("H132").Value = 14357985
$ws.Range("I132").Value = 14357985
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 43073955
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -43071425
$ws.Range("H137").Value = 29909.227
$ws.Range("I137").Value = 51906.668
$ws.Range("K137").Value = 155720.004
$ws.Range("M137").Value = -153170.004
$ws.Range("H138").Value = 25401.217
$ws.Range("I138").Value = 3113.762
$ws.Range("J138").Value = 44122.68
$ws.Range("K138").Value = 9341.286
$ws.Range("L138").Value = 132368.04
$ws.Range("M138").Value = -4201.286
$ws.Range("N138").Value = -142648.04

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H45").Value = 6999.3335
$ws.Range("I45").Value = 1198
$ws.Range("J45").Value = 9900
$ws.Range("K45").Value = 1198
$ws.Range("L45").Value = 9900
$ws.Range("M45").Value = -821
$ws.Range("N45").Value = -10654
$ws.Range("H61").Value = 8666.044
$ws.Range("I61").Value = 4453.6313
$ws.Range("K61").Value = 4453.6313
$ws.Range("M61").Value = -4241.6313
$ws.Range("H132").Value = 2842.7568
$ws.Range("I132").Value = 2640.4546
$ws.Range("J132").Value = 3139.4666
$ws.Range("K132").Value = 7921.3638
$ws.Range("L132").Value = 9418.399800000001
$ws.Range("M132").Value = -5391.3638
$ws.Range("N132").Value = -14478.3998
$ws.Range("H136").Value = 8666.044
$ws.Range("I136").Value = 4453.6313
$ws.Range("K136").Value = 13360.8939
$ws.Range("M136").Value = -10810.8939

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3125790
$ws.Range("I94").Value = 761.8333
$ws.Range("K94").Value = 761.8333
$ws.Range("M94").Value = -310.8333
$ws.Range("H105").Value = 3700.85
$ws.Range("I105").Value = 2900
$ws.Range("K105").Value = 2900
$ws.Range("M105").Value = -1153
$ws.Range("H132").Value = 69999
$ws.Range("J132").Value = 69999
$ws.Range("L132").Value = 69999
$ws.Range("N132").Value = -80119
$ws.Range("H134").Value = 2978.9592
$ws.Range("I134").Value = 3037.0232
$ws.Range("K134").Value = 9111.069600000001
$ws.Range("M134").Value = -6576.069600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 127038.375
$ws.Range("I132").Value = 144472.42
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 433417.26
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -430887.26
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 149.33333
$ws.Range("I14").Value = 149.33333
$ws.Range("K14").Value = 447.99999
$ws.Range("M14").Value = -274.99999
$ws.Range("H34").Value = 5422.846
$ws.Range("I34").Value = 1200
$ws.Range("J34").Value = 8062.125
$ws.Range("K34").Value = 3600
$ws.Range("L34").Value = 24186.375
$ws.Range("M34").Value = -3516
$ws.Range("N34").Value = -24354.375
$ws.Range("H39").Value = 11972.875
$ws.Range("I39").Value = 12826.143
$ws.Range("J39").Value = 6000
$ws.Range("K39").Value = 38478.429
$ws.Range("L39").Value = 18000
$ws.Range("M39").Value = -38184.429
$ws.Range("N39").Value = -18588
$ws.Range("H55").Value = 2752.5
$ws.Range("I55").Value = 2003
$ws.Range("J55").Value = 5001
$ws.Range("K55").Value = 6009
$ws.Range("L55").Value = 15003
$ws.Range("M55").Value = -5832
$ws.Range("N55").Value = -15357
$ws.Range("H92").Value = 2045.9231
$ws.Range("J92").Value = 1460
$ws.Range("L92").Value = 4380
$ws.Range("N92").Value = -6876
$ws.Range("H122").Value = 1279.8
$ws.Range("I122").Value = 499.6
$ws.Range("J122").Value = 2060
$ws.Range("K122").Value = 4496.400000000001
$ws.Range("L122").Value = 18540
$ws.Range("M122").Value = -2046.400000000001
$ws.Range("N122").Value = -23440

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 426.87878
$ws.Range("I2").Value = 543.6667
$ws.Range("J2").Value = 286.73334
$ws.Range("K2").Value = 543.6667
$ws.Range("L2").Value = 286.73334
$ws.Range("M2").Value = -430.6667
$ws.Range("N2").Value = -512.73334
$ws.Range("H80").Value = 3954.5454
$ws.Range("I80").Value = 3139.077
$ws.Range("J80").Value = 5132.4443
$ws.Range("K80").Value = 3139.077
$ws.Range("L80").Value = 5132.4443
$ws.Range("M80").Value = -2141.077
$ws.Range("N80").Value = -7128.4443
$ws.Range("H83").Value = 3954.5454
$ws.Range("I83").Value = 3139.077
$ws.Range("J83").Value = 5132.4443
$ws.Range("K83").Value = 15695.385
$ws.Range("L83").Value = 25662.2215
$ws.Range("M83").Value = -10703.385
$ws.Range("N83").Value = -35646.2215
$ws.Range("H126").Value = 2739.2856
$ws.Range("I126").Value = 2100
$ws.Range("K126").Value = 6300
$ws.Range("M126").Value = -3830
$ws.Range("H134").Value = 65518.75
$ws.Range("J134").Value = 65518.75
$ws.Range("L134").Value = 196556.25
$ws.Range("N134").Value = -201626.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1863.4688
$ws.Range("J55").Value = 2615.1333
$ws.Range("L55").Value = 2615.1333
$ws.Range("N55").Value = -2961.1333
$ws.Range("H68").Value = 2177.8572
$ws.Range("I68").Value = 1891
$ws.Range("J68").Value = 2895
$ws.Range("K68").Value = 1891
$ws.Range("L68").Value = 2895
$ws.Range("M68").Value = -1142
$ws.Range("N68").Value = -4393
$ws.Range("H71").Value = 2177.8572
$ws.Range("I71").Value = 1891
$ws.Range("J71").Value = 2895
$ws.Range("K71").Value = 9455
$ws.Range("L71").Value = 14475
$ws.Range("M71").Value = -5711
$ws.Range("N71").Value = -21963
$ws.Range("H122").Value = 4407.375
$ws.Range("I122").Value = 4089.75
$ws.Range("J122").Value = 4725
$ws.Range("K122").Value = 12269.25
$ws.Range("L122").Value = 14175
$ws.Range("M122").Value = -9819.25
$ws.Range("N122").Value = -19075
$ws.Range("H132").Value = 6113.324
$ws.Range("I132").Value = 5799.6
$ws.Range("K132").Value = 17398.8
$ws.Range("M132").Value = -14868.8
$ws.Range("H136").Value = 5799.5713
$ws.Range("I136").Value = 5099.6665
$ws.Range("K136").Value = 15298.9995
$ws.Range("M136").Value = -12748.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 26195.742
$ws.Range("I122").Value = 29174.686
$ws.Range("K122").Value = 87524.058
$ws.Range("M122").Value = -85074.058
$ws.Range("H136").Value = 12480.972
$ws.Range("I136").Value = 13065.141
$ws.Range("J136").Value = 7140
$ws.Range("K136").Value = 39195.423
$ws.Range("L136").Value = 21420
$ws.Range("M136").Value = -36645.423
$ws.Range("N136").Value = -26520
